$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D10 ("0.650") would otherwise be auto-detected as the number 0.65 and
# lose its trailing zero, so pin it to text before assigning the value.
$ws.Range("D10").NumberFormat = "@"

$ws.Range("D2").Value = '67.132.72'
$ws.Range("E2").Value = '  -1.77%  '
$ws.Range("D3").Value = '3.490.63'
$ws.Range("E3").Value = '  -3.60%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '197.65'
$ws.Range("E5").Value = '  +1.52%  '
$ws.Range("D6").Value = '543.36'
$ws.Range("E6").Value = '  -7.15%  '
$ws.Range("D7").Value = '3.484.40'
$ws.Range("E7").Value = '  -3.64%  '
$ws.Range("D8").Value = '0.601'
$ws.Range("E8").Value = '  -3.35%  '
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("D10").Value = '0.650'
$ws.Range("E10").Value = '  -4.52%  '
$ws.Range("D11").Value = '62.92'
$ws.Range("E11").Value = '  +12.74%  '
$ws.Range("D12").Value = '0.141'
$ws.Range("E12").Value = '  -7.08%  '
$ws.Range("D13").Value = '0.0000268'
$ws.Range("E13").Value = '  -8.05%  '
$ws.Range("D14").Value = '9.71'
$ws.Range("E14").Value = '  -3.62%  '
$ws.Range("D15").Value = '4.048.02'
$ws.Range("E15").Value = '  -3.39%  '
$ws.Range("D16").Value = '3.485.28'
$ws.Range("E16").Value = '  -3.65%  '
$ws.Range("D17").Value = '0.123'
$ws.Range("E17").Value = '  -1.83%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '18.32'
$ws.Range("E18").Value = '  -1.41%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '66.758.93'
$ws.Range("E19").Value = '  -2.11%  '
$ws.Range("D20").Value = '11.72'
$ws.Range("E20").Value = '  -6.58%  '
$ws.Range("E21").Value = '  -5.51%  '
$ws.Range("D22").Value = '388.99'
$ws.Range("E22").Value = '  -3.86%  '
$ws.Range("E23").Value = '  -6.89%  '
$ws.Range("D24").Value = '11.76'
$ws.Range("E24").Value = '  -11.81%  '
$ws.Range("D25").Value = '81.86'
$ws.Range("E25").Value = '  -5.07%  '
$ws.Range("D26").Value = '2.78'
$ws.Range("E26").Value = '  -6.06%  '
$ws.Range("D27").Value = '12.04'
$ws.Range("D28").Value = '3.71'
$ws.Range("E28").Value = '  -7.57%  '
$ws.Range("D29").Value = '8.71'
$ws.Range("E29").Value = '  -5.24%  '
$ws.Range("D30").Value = '30.67'
$ws.Range("E30").Value = '  -3.31%  '
$ws.Range("D31").Value = '676.52'
$ws.Range("E31").Value = '  -2.40%  '
$ws.Range("E32").Value = '  -13.79%  '
$ws.Range("D33").Value = '11.61'
$ws.Range("E33").Value = '  -5.50%  '
$ws.Range("D34").Value = '63.18'
$ws.Range("E35").Value = '  -7.45%  '
$ws.Range("D36").Value = '38.52'
$ws.Range("E36").Value = '  -10.04%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D38").Value = '0.397'
$ws.Range("E38").Value = '  -5.19%  '
$ws.Range("D39").Value = '0.131'
$ws.Range("E39").Value = '  -3.47%  '
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("D41").Value = '3.054.51'
$ws.Range("E41").Value = '  -3.36%  '
$ws.Range("D43").Value = '0.0₃0670'
$ws.Range("E43").Value = '  -15.63%  '
$ws.Range("E44").Value = '  -14.64%  '
$ws.Range("D45").Value = '2.72'
$ws.Range("E45").Value = '  +4.39%  '
$ws.Range("D46").Value = '2.66'
$ws.Range("E46").Value = '  +2.60%  '
$ws.Range("E47").Value = '  -7.50%  '
$ws.Range("E48").Value = '  -4.78%  '
$ws.Range("D49").Value = '137.44'
$ws.Range("E49").Value = '  -4.34%  '
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").Value = '8.13'
$ws.Range("E50").Value = '  -8.48%  '
$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D51").Value = '2.87'
$ws.Range("E51").Value = '  -7.86%  '
